$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

$ws.Range('AC9').Value = 'GND,Net-(U1-AOUT)'
$ws.Range('AD9').Value = 'GND,Net-(U1-AOUT)'
$ws.Range('AC11').Value = 'GND,+5V'
$ws.Range('AD11').Value = 'GND,+5V'
$ws.Range('AE11').Value = 'vcc,Default'
$ws.Range('AC13').Value = 'GND,+5V'
$ws.Range('AD13').Value = 'GND,+5V'
$ws.Range('AE13').Value = 'vcc,Default'
$ws.Range('AC15').Value = 'Net-(JP1-Pin_3),Net-(JP1-Pin_1),Net-(JP1-Pin_7),Net-(JP1-Pin_5),Net-(JP1-Pin_11),+5V,Net-(JP1-Pin_9)'
$ws.Range('AD15').Value = 'Net-(JP1-Pin_3),Net-(JP1-Pin_1),Net-(JP1-Pin_7),Net-(JP1-Pin_5),Net-(JP1-Pin_11),+5V,Net-(JP1-Pin_9)'
$ws.Range('AE15').Value = 'vcc,Default'
$ws.Range('AC16').Value = 'unconnected-(P1-Pin_36-Pad36),/A6,/A7,/~{WR},/A4,/D4,/D0,/D1,unconnected-(P1-Pin_23-Pad23),GND,/A2,unconnected-(P1-Pin_3-Pad3),/D5,unconnected-(P1-Pin_4-Pad4),unconnected-(P1-Pin_35-Pad35),/D3,unconnected-(P1-Pin_38-Pad38),/D2,unconnected-(P1-Pin_8-Pad8),/~{RD},/D6,unconnected-(P1-Pin_39-Pad39),unconnected-(P1-Pin_1-Pad1),unconnected-(P1-Pin_21-Pad21),unconnected-(P1-Pin_37-Pad37),unconnected-(P1-Pin_2-Pad2),unconnected-(P1-Pin_7-Pad7),unconnected-(P1-Pin_19-Pad19),/D7,/A0,unconnected-(P1-Pin_5-Pad5),/A1,unconnected-(P1-Pin_6-Pad6),/A3,/~{INT},/A5,/~{IORQ},+5V,/~{RESET}'
$ws.Range('AD16').Value = '~{RESET}'
$ws.Range('AE16').Value = 'vcc,Default'
$ws.Range('AC19').Value = 'GND,Net-(CON1-PadT)'
$ws.Range('AD19').Value = 'GND,Net-(CON1-PadT)'
$ws.Range('AC20').Value = 'GND,Net-(JP1-Pin_3),/A2,/A3,Net-(JP1-Pin_1),Net-(JP1-Pin_7),/A6,/A7,Net-(JP1-Pin_5),/A5,/A4,Net-(JP1-Pin_11),/~{IORQ},+5V,Net-(JP1-Pin_9),/~{CS}'
$ws.Range('AD20').Value = '~{CS}'
$ws.Range('AE20').Value = 'vcc,Default'
$ws.Range('AC21').Value = 'GND,/AUDIO_CH2,/L,Net-(U1-MP),Net-(U2B--),Net-(U1-CV),/AUDIO_CH1,Net-(U1-AOUT),+5V,/R'
$ws.Range('AD21').Value = 'R'
$ws.Range('AE21').Value = 'vcc,Default'
$ws.Range('AC22').Value = 'GND,/AUDIO_CH2,Net-(U1-MP),/AUDIO_CH1,Net-(U1-CV),/SMPAC,/DAC_CLK,Net-(U1-SWIN),Net-(U1-AOUT),/DOAB,+5V,/SMPBD,unconnected-(U1-TST2-Pad15)'
$ws.Range('AD22').Value = 'SMPBD,unconnected-(U1-TST2-Pad15)'
$ws.Range('AE22').Value = 'vcc,Default'
$ws.Range('AC23').Value = '/~{WR},/D4,/SMPAC,/D0,/D1,/~{CS},GND,/D5,/D3,/DAC_CLK,/DOAB,/D2,/~{RD},/D6,/D7,/A0,unconnected-(U4-TEST-Pad9),/SMPBD,/A1,unconnected-(U4-DOCD-Pad22),+5V,unconnected-(U4-~{IRQ}-Pad2),/CLK,/~{RESET}'
$ws.Range('AD23').Value = '~{RESET}'
$ws.Range('AE23').Value = 'vcc,Default'
$ws.Range('AC24').Value = 'GND,/CLK,+5V'
$ws.Range('AD24').Value = 'CLK,+5V'
$ws.Range('AE24').Value = 'vcc,Default'

$ws.Rows.Item(16).RowHeight = 180
